$wb = $excel.ActiveWorkbook

# --- Sheet "保險" (insurance, worksheet #5) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B1").Value = 'company'
$ws5.Range("C1").Value = 'name'
$ws5.Range("D1").Value = 'owner'
$ws5.Range("E1").Value = 'property_category'
$ws5.Range("F1").Value = 'category'
$ws5.Range("G1").Value = 'date'
$ws5.Range("H1").Value = 'legislator_name'
$ws5.Range("I1").Value = 'legislator_id'
$ws5.Range("J1").Value = 'source_file'
$ws5.Range("K1").Value = 'index'
$ws5.Range("A2").Value = 105
$ws5.Range("B2").Value = '三商美邦人壽'
$ws5.Range("C2").Value = '世紀理財變額萬能終身壽險'
$ws5.Range("D2").Value = '蔡煌瑯'
$ws5.Range("E2").Value = 'insurance'
$ws5.Range("F2").Value = 'normal'
$ws5.Range("G2").Value = '2012-04-27'
$ws5.Range("H2").Value = '蔡煌瑯'
$ws5.Range("I2").Value = 752
$ws5.Range("J2").Value = 'tmpd4981'
$ws5.Range("K2").Value = 105
$ws5.Range("A3").Value = 106
$ws5.Range("B3").Value = '三商美邦人壽'
$ws5.Range("C3").Value = '世紀理財變額萬能終身壽險'
$ws5.Range("D3").Value = '王琴賀'
$ws5.Range("E3").Value = 'insurance'
$ws5.Range("F3").Value = 'normal'
$ws5.Range("G3").Value = '2012-04-27'
$ws5.Range("H3").Value = '蔡煌瑯'
$ws5.Range("I3").Value = 752
$ws5.Range("J3").Value = 'tmpd4981'
$ws5.Range("K3").Value = 106
$ws5.Range("A4").Value = 107
$ws5.Range("B4").Value = '新光人壽'
$ws5.Range("C4").Value = '美利外幣終生還本型保險'
$ws5.Range("D4").Value = '王琴賀'
$ws5.Range("E4").Value = 'insurance'
$ws5.Range("F4").Value = 'normal'
$ws5.Range("G4").Value = '2012-04-27'
$ws5.Range("H4").Value = '蔡煌瑯'
$ws5.Range("I4").Value = 752
$ws5.Range("J4").Value = 'tmpd4981'
$ws5.Range("K4").Value = 107

# --- Sheet "債務" (debt, worksheet #6) ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B1").Value = 'species'
$ws6.Range("C1").Value = 'debtor'
$ws6.Range("D1").Value = 'owner'
$ws6.Range("E1").Value = 'total'
$ws6.Range("F1").Value = 'register_date'
$ws6.Range("G1").Value = 'register_reason'
$ws6.Range("H1").Value = 'property_category'
$ws6.Range("I1").Value = 'category'
$ws6.Range("J1").Value = 'date'
$ws6.Range("K1").Value = 'legislator_name'
$ws6.Range("L1").Value = 'legislator_id'
$ws6.Range("M1").Value = 'source_file'
$ws6.Range("N1").Value = 'index'
$ws6.Range("A2").Value = 117
$ws6.Range("B2").Value = '房屋貸款'
$ws6.Range("C2").Value = '王琴賀'
$ws6.Range("D2").Value = '台中商業銀行埔里分行南投縣埔里鎮西康路'
$ws6.Range("E2").Value = 4054661
$ws6.Range("F2").Value = '95年01月23日'
$ws6.Range("G2").Value = '設定'
$ws6.Range("H2").Value = 'debt'
$ws6.Range("I2").Value = 'normal'
$ws6.Range("J2").Value = '2012-04-27'
$ws6.Range("K2").Value = '蔡煌瑯'
$ws6.Range("L2").Value = 752
$ws6.Range("M2").Value = 'tmpd4981'
$ws6.Range("N2").Value = 117
$ws6.Range("A3").Value = 118
$ws6.Range("B3").Value = '房屋貸款'
$ws6.Range("C3").Value = '王琴賀'
$ws6.Range("D3").Value = '台中商業銀行埔里分行南投縣埔里鎮西康路'
$ws6.Range("E3").Value = 344914
$ws6.Range("F3").Value = '98年11月16日'
$ws6.Range("G3").Value = '設定'
$ws6.Range("H3").Value = 'debt'
$ws6.Range("I3").Value = 'normal'
$ws6.Range("J3").Value = '2012-04-27'
$ws6.Range("K3").Value = '蔡煌瑯'
$ws6.Range("L3").Value = 752
$ws6.Range("M3").Value = 'tmpd4981'
$ws6.Range("N3").Value = 118
$ws6.Range("A4").Value = 119
$ws6.Range("B4").Value = '房屋貸款'
$ws6.Range("C4").Value = '蔡煌瑯'
$ws6.Range("D4").Value = '臺灣銀行臺北市中正區重慶南路'
$ws6.Range("E4").Value = 14193131
$ws6.Range("F4").Value = '98年01月16日'
$ws6.Range("G4").Value = '設定'
$ws6.Range("H4").Value = 'debt'
$ws6.Range("I4").Value = 'normal'
$ws6.Range("J4").Value = '2012-04-27'
$ws6.Range("K4").Value = '蔡煌瑯'
$ws6.Range("L4").Value = 752
$ws6.Range("M4").Value = 'tmpd4981'
$ws6.Range("N4").Value = 119
$ws6.Range("A5").Value = 121
$ws6.Range("B5").Value = '房屋貸款'
$ws6.Range("C5").Value = '王琴賀'
$ws6.Range("D5").Value = '台中商業銀行埔里分行南投縣埔里鎮西康路'
$ws6.Range("E5").Value = 955272
$ws6.Range("F5").Value = '99年01月15日'
$ws6.Range("G5").Value = '設定'
$ws6.Range("H5").Value = 'debt'
$ws6.Range("I5").Value = 'normal'
$ws6.Range("J5").Value = '2012-04-27'
$ws6.Range("K5").Value = '蔡煌瑯'
$ws6.Range("L5").Value = 752
$ws6.Range("M5").Value = 'tmpd4981'
$ws6.Range("N5").Value = 121
$ws6.Range("A6").Value = 122
$ws6.Range("B6").Value = '信用貸款'
$ws6.Range("C6").Value = '王琴賀'
$ws6.Range("D6").Value = '台中商業銀行埔里分行南投縣埔里鎮西康路'
$ws6.Range("E6").Value = 3000000
$ws6.Range("F6").Value = '100年01月19曰'
$ws6.Range("G6").Value = '設定'
$ws6.Range("H6").Value = 'debt'
$ws6.Range("I6").Value = 'normal'
$ws6.Range("J6").Value = '2012-04-27'
$ws6.Range("K6").Value = '蔡煌瑯'
$ws6.Range("L6").Value = 752
$ws6.Range("M6").Value = 'tmpd4981'
$ws6.Range("N6").Value = 122
